$d = $word.ActiveDocument

# 1. Standard software paragraph
$d.Content.Find.Execute(
    "the instance user, is super-user, so can install whatever else is needed using",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "the instance user has super-user rights, so you can install whatever else is needed using",
    2)

# 2. Genomics data paragraph
$d.Content.Find.Execute(
    "In the ~/data folder genomic data from 51 example-genomes are found in standard impute-me format: The",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "In the ~/data folder. Genomic data from 51 example-genomes are available in standard impute-me format: The",
    2)

# 3. Art repository data paragraph
# (Find text is deliberately started one character past the run boundary -
#  right after the leading colon - so the replacement does not inherit the
#  bold formatting of the preceding "Art repository data" run.)
$d.Content.Find.Execute(
    "Large image repositories are available online, and finding the right one may be key to this project. However, two large sets are already pre-loaded on the computer: The",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "In the ~/art folder. Large image repositories are available online, and finding the right one may be key to this project. However, two large sets are already pre-loaded on the computer: The",
    2)

Write-Host "Done simple replacements"

# 4. Custom software paragraph - larger restructuring plus hyperlink move.
# 4a. Remove the old "from here" hyperlink and merge surrounding text:
#     "running code (from here), the ~/programs" -> "running code for impute.me, the ~/programs"
$ok4a = $d.Content.Find.Execute(
    "running code (from here), the ~/programs",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "running code for impute.me, the ~/programs",
    2)
Write-Host "4a: $ok4a"

# 4b. Add "In the ~/srv and ~/programs folder." sentence after "Custom software:"
# (again start the Find text one character past the run boundary so the new
#  text does not pick up the bold formatting of "Custom software")
$ok4b = $d.Content.Find.Execute(
    "Each computer has a pre-configured copy of the impute.me web-analysis server running on shiny.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "In the ~/srv and ~/programs folder. Each computer has a pre-configured copy of the impute.me web-analysis server running on shiny.",
    2)
Write-Host "4b: $ok4b"

# 4c. Replace the trailing ")" (the lone plain run right after the verbatim URL)
#     with the new closing sentence, leaving the verbatim-styled URL run untouched.
$rng4c = $d.Content
$ok4c = $rng4c.Find.Execute(
    "amazonaws.com:3838",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "4c found: $ok4c"
$rng4c.Collapse(0)
$rng4c.MoveEnd(1, 1)
$rng4c.Text = "). To work with this code, it is recommended to make a fork of"

# 4d. Insert a single space as its own run after "make a fork of"
$rng4d = $d.Content
$ok4d = $rng4d.Find.Execute("make a fork of", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "4d found: $ok4d"
$rng4d.Collapse(0)
$rng4d.InsertAfter(" ")

# 4e. Insert placeholder text (to become the new hyperlink) plus the trailing
#     sentence, then convert the placeholder into a hyperlink re-using the
#     same target as the old "from here" link.
$rng4e = $d.Content
$ok4e = $rng4e.Find.Execute("make a fork of ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "4e found: $ok4e"
$rng4e.Collapse(0)
$rng4e.InsertAfter("GITHUBREPOPLACEHOLDER and put locally.")

$rng4f = $d.Content
$ok4f = $rng4f.Find.Execute("GITHUBREPOPLACEHOLDER", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "4f found: $ok4f"
$d.Hyperlinks.Add($rng4f, "https://github.com/lassefolkersen/impute-me", "", "", "the github repository") | Out-Null

Write-Host "Done custom software paragraph"
